$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (test_number) for rows 2-11 from 2 to 4 (sample 4)
$ws.Range("A2:A11").Value = 4

# q7 (row 8) answers_4 changes from 5 to 4
$ws.Range("E8").Value = 4

# q10 (row 11) answers_4 changes from 1 to 6
$ws.Range("E11").Value = 6

# Update the selection to match the edited range
$ws.Range("A2:A11").Select()
